$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('AK_1.png', 0.545, 0.002, 0.982, 0.333, 'Tidak Diketahui', 'Salah'),
    @('AK_2.png', 0.9370000000000001, 0.003, 0.988, 0.4, 'Tidak Diketahui', 'Salah'),
    @('AK_3.png', 0.901, 0.003, 0.979, 0.8, 'Akhlak Kamiswara', 'Benar'),
    @('AK_4.png', 0.552, 0.002, 0.976, 0.533, 'Akhlak Kamiswara', 'Benar'),
    @('AK_5.png', 0.5629999999999999, 0.002, 0.99, 0.333, 'Tidak Diketahui', 'Salah'),
    @('MIB_1.png', 1.257, 0.004, 0.963, 0.467, 'Tidak Diketahui', 'Salah'),
    @('MIB_2.png', 1.098, 0.004, 0.972, 0.467, 'Tidak Diketahui', 'Salah'),
    @('MIB_3.png', 1.354, 0.004, 0.999, 1, 'Muhammad Iqbal Baqi', 'Benar'),
    @('MIB_4.png', 1.113, 0.004, 0.984, 0.733, 'Muhammad Iqbal Baqi', 'Benar'),
    @('MIB_5.png', 1.262, 0.004, 0.98, 0.733, 'Muhammad Iqbal Baqi', 'Benar'),
    @('AAH_1.png', 0.892, 0.003, 0.974, 0.6, 'Andrea Ayunove Hutami', 'Benar'),
    @('AAH_2.png', 1.097, 0.004, 0.992, 0.667, 'Andrea Ayunove Hutami', 'Benar'),
    @('AAH_3.png', 0.862, 0.003, 0.975, 0.6, 'Andrea Ayunove Hutami', 'Benar'),
    @('TI_1.png', 0.953, 0.003, 0.949, 0.533, 'Toni Ismail', 'Benar'),
    @('TI_2.png', 0.944, 0.003, 0.989, 0.4, 'Tidak Diketahui', 'Salah'),
    @('TI_3.png', 0.785, 0.003, 0.991, 0.6, 'Toni Ismail', 'Benar'),
    @('TI_4.png', 0.753, 0.003, 0.974, 0.733, 'Toni Ismail', 'Benar'),
    @('TI_5.png', 1.064, 0.004, 0.991, 0.533, 'Toni Ismail', 'Benar'),
    @('RAS_1.png', 0.712, 0.002, 0.975, 0.333, 'Tidak Diketahui', 'Salah'),
    @('RAS_2.png', 1.029, 0.003, 0.969, 0.6, 'Ridha Ayu Salsabila', 'Benar'),
    @('RAS_3.png', 0.5629999999999999, 0.002, 0.982, 0.4, 'Tidak Diketahui', 'Salah'),
    @('RAS_4.png', 1.136, 0.004, 0.967, 0.267, 'Tidak Diketahui', 'Salah'),
    @('RAS_5.png', 1.094, 0.004, 0.978, 0.733, 'Ridha Ayu Salsabila', 'Benar'),
    @('RR_1.png', 1.157, 0.004, 0.987, 0.6, 'Rafiqo Rapitasari', 'Benar'),
    @('RR_2.png', 1.193, 0.004, 0.988, 0.533, 'Rafiqo Rapitasari', 'Benar'),
    @('RR_3.png', 0.996, 0.003, 0.9, 0.4, 'Tidak Diketahui', 'Salah'),
    @('RR_4.png', 1.19, 0.004, 0.992, 0.667, 'Rafiqo Rapitasari', 'Benar'),
    @('RR_5.png', 1.206, 0.004, 0.988, 0.8, 'Rafiqo Rapitasari', 'Benar'),
    @('AR_1.png', 0.891, 0.003, 0.965, 0.733, 'Arizli Romadhon', 'Benar'),
    @('GA_1.png', 1.207, 0.004, 0.979, 1, 'Gege Ardiyansyah', 'Benar'),
    @('GA_2.png', 0.6840000000000001, 0.002, 0.99, 0.267, 'Tidak Diketahui', 'Salah'),
    @('GA_3.png', 0.758, 0.003, 0.99, 0.333, 'Tidak Diketahui', 'Salah'),
    @('FY_1.png', 1.055, 0.004, 0.976, 0.333, 'Tidak Diketahui', 'Salah'),
    @('FY_2.png', 1.269, 0.004, 0.958, 0.333, 'Tidak Diketahui', 'Salah'),
    @('FY_3.png', 1.238, 0.004, 0.992, 0.333, 'Tidak Diketahui', 'Salah'),
    @('FY_4.png', 1.163, 0.004, 0.979, 0.333, 'Tidak Diketahui', 'Salah'),
    @('TO_1.png', 0.902, 0.003, 0.971, 0.4, 'Tidak Diketahui', 'Salah'),
    @('TO_2.png', 0.9399999999999999, 0.003, 0.986, 0.6, 'Tiara Oktavian', 'Benar'),
    @('TO_3.png', 0.833, 0.003, 0.983, 0.6, 'Tiara Oktavian', 'Benar'),
    @('TO_4.png', 2.36, 0.008, 0.799, 0.333, 'Tidak Diketahui', 'Salah'),
    @('TO_5.png', 2.58, 0.008, 0.78, 0.467, 'Tidak Diketahui', 'Salah'),
    @('TD_1.png', 1.565, 0.005, 0.826, 0.4, 'Tidak Diketahui', 'Benar'),
    @('TD_2.png', 1.574, 0.005, 0.763, 0.467, 'Tidak Diketahui', 'Benar'),
    @('TD_3.png', 1.035, 0.003, 0.971, 0.2, 'Tidak Diketahui', 'Benar'),
    @('TD_4.png', 0.977, 0.003, 0.988, 0.267, 'Tidak Diketahui', 'Benar'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}
